# Erstanmeldung in Sametime.pptx - apply 2024-07-24 update
#
# 1) The cached "today" text of the datetimeFigureOut date-placeholder
#    field is refreshed from 10.07.2024 to 24.07.2024 on the slide
#    master and on every slide layout (12 placeholders total).
# 2) A new text run "24-07.2024" is added to the (until now empty)
#    subtitle placeholder on slide 1.

$p = $ppt.ActivePresentation

function Update-DateText($shape, [string]$newText) {
    if ($shape.Type -eq 14) {
        $phf = $shape.PlaceholderFormat
        if ($phf -ne $null -and $phf.Type -eq 16) {
            $shape.TextFrame.TextRange.Text = $newText
        }
    }
}

# -- Slide master date placeholder --
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateText $master.Shapes.Item($i) "24.07.2024"
}

# -- Every slide layout's date placeholder --
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateText $layout.Shapes.Item($i) "24.07.2024"
    }
}

# -- Slide 1: add "24-07.2024" run to the empty subtitle placeholder --
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$newText = "24-07.2024"
[void]$subtitle.TextFrame.TextRange.InsertBefore($newText)
$newRun = $subtitle.TextFrame.TextRange.Characters(1, $newText.Length)
$newRun.LanguageID = "de-DE"
